$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4149.185064510225
$ws.Range("C3").Value = 4149.185064510225
$ws.Range("C4").Value = 3981.336150629966
$ws.Range("C5").Value = 3981.336150629966
$ws.Range("C6").Value = 3981.336150629966
$ws.Range("C7").Value = 3717.908737070632
$ws.Range("C8").Value = 3717.908737070632
$ws.Range("C9").Value = 3707.438063422706
$ws.Range("C10").Value = 3707.438063422706
$ws.Range("C11").Value = 3707.438063422706
$ws.Range("C12").Value = 3659.841069047226
